# Generate Report for Handback
# Updates the Overview / zh-cn / de-de sheets to reflect that both
# language handbacks are now complete:
#   - Status text switches from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears.
#   - The zh-cn and de-de sheets gain a "Latest Target File" hyperlink
#     (col I) pointing at the source .md file, a "Latest Handback File"
#     name (col J) and an updated "Latest Handback DateTime" (col K).

$wb = $excel.ActiveWorkbook

$srcMdName = "f76c3241-b22a-44c1-95af-8ae9c78f4368.md"
$srcMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ffbaf34739bb26b21fe8201f858498e908efb0a6/e2e/f76c3241-b22a-44c1-95af-8ae9c78f4368.md"
$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---- Status column updates (shared text, ripples to every cell that
#      currently shows "Ready for handoff") --------------------------------
$wsOverview.Cells.Item(2, 5).Value2 = $statusText   # Overview!E2 (zh-cn status)
$wsOverview.Cells.Item(2, 6).Value2 = $statusText   # Overview!F2 (de-de status)
$wsZhCn.Cells.Item(2, 3).Value2 = $statusText       # zh-cn!C2 (Status)
$wsDeDe.Cells.Item(2, 3).Value2 = $statusText       # de-de!C2 (Status)

# ---- zh-cn row 2: Latest Target File / Latest Handback File / DateTime ---
$wsZhCn.Cells.Item(2, 9).Value2 = $srcMdName
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(2, 9), $srcMdUrl, "", "", $srcMdName) | Out-Null
$wsZhCn.Cells.Item(2, 9).Font.Underline = 2
$wsZhCn.Cells.Item(2, 9).Font.Color = 15570276

$wsZhCn.Cells.Item(2, 10).Value2 = "f76c3241-b22a-44c1-95af-8ae9c78f4368.c3a57fbd94d88c55539ed187636a77d6a0ded828.zh-cn.xlf"
$wsZhCn.Cells.Item(2, 11).Value2 = "2016-08-31 13:02:17"

# ---- de-de row 2: Latest Target File / Latest Handback File / DateTime ---
$wsDeDe.Cells.Item(2, 9).Value2 = $srcMdName
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(2, 9), $srcMdUrl, "", "", $srcMdName) | Out-Null
$wsDeDe.Cells.Item(2, 9).Font.Underline = 2
$wsDeDe.Cells.Item(2, 9).Font.Color = 15570276

$wsDeDe.Cells.Item(2, 10).Value2 = "f76c3241-b22a-44c1-95af-8ae9c78f4368.c3a57fbd94d88c55539ed187636a77d6a0ded828.de-de.xlf"
$wsDeDe.Cells.Item(2, 11).Value2 = "2016-08-31 13:02:38"

# ---- Column width adjustments (grow to fit the now-longer text) ----------
# Target widths land on this engine's internal 1/6-character grid, so we
# dial the ColumnWidth input back by the fixed 5/6 padding it re-adds.
$wsOverview.Columns.Item(5).ColumnWidth = 29.16666667   # -> ~29.98 (Status, zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666667   # -> ~29.98 (Status, de-de)

$wsZhCn.Columns.Item(3).ColumnWidth = 29.16666667        # -> ~29.98 (Status)
$wsZhCn.Columns.Item(9).ColumnWidth = 39.16666667        # -> 40 (Latest Target File)
$wsZhCn.Columns.Item(10).ColumnWidth = 39.16666667       # -> 40 (Latest Handback File)

$wsDeDe.Columns.Item(3).ColumnWidth = 29.16666667        # -> ~29.98 (Status)
$wsDeDe.Columns.Item(9).ColumnWidth = 39.16666667        # -> 40 (Latest Target File)
$wsDeDe.Columns.Item(10).ColumnWidth = 39.16666667       # -> 40 (Latest Handback File)
